# Daily Satellite Data Update
# Adds the 30.12.2025 flyover rows, refreshes the cloud-coverage figures for
# the existing rows, and extends the conditional-formatting ranges to cover
# the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Refresh cloud-coverage (O:R) figures for the existing 27/28/29 Dec rows
# ---------------------------------------------------------------------
$ws.Range("O3").Value = 94
$ws.Range("P3").Value = 66
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 26

$ws.Range("O4").Value = 67
$ws.Range("P4").Value = 65
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0

$ws.Range("O5").Value = 84
$ws.Range("P5").Value = 84
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0

# Re-colour the cloud-coverage cells to match the refreshed figures.
$ws.Range("O3").Interior.Color = 8351984   # F0707F
$ws.Range("P3").Interior.Color = 15521732  # C4D7EC
$ws.Range("Q3").Interior.Color = 13143125  # 558CC8
$ws.Range("R3").Interior.Color = 14068096  # 80A9D6

$ws.Range("O4").Interior.Color = 8351984   # F0707F
$ws.Range("P4").Interior.Color = 15521732  # C4D7EC
$ws.Range("Q4").Interior.Color = 13143125  # 558CC8
$ws.Range("R4").Interior.Color = 13143125  # 558CC8

$ws.Range("O5").Interior.Color = 8351984   # F0707F
$ws.Range("P5").Interior.Color = 16248550  # E6EEF7
$ws.Range("Q5").Interior.Color = 13143125  # 558CC8
$ws.Range("R5").Interior.Color = 13143125  # 558CC8

# ---------------------------------------------------------------------
# 2. Append the two new flyovers for 30.12.2025
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "30.12.2025"
$ws.Range("B6").Value = 18
$ws.Range("C6").Value = "05:03"
$ws.Range("D6").Value = "10:19"
$ws.Range("E6").Value = "04:30:53"
$ws.Range("F6").Value = "04:33:31"
$ws.Range("G6").Value = "04:36:02"
$ws.Range("H6").Value = "04:38:34"
$ws.Range("I6").Value = "04:41:12"
$ws.Range("J6").Value = "15°"
$ws.Range("K6").Value = "04:34:45"
$ws.Range("L6").Value = -22.5
$ws.Range("M6").Value = "A+B"
$ws.Range("N6").Value = "2"

$ws.Range("A7").Value = "30.12.2025"
$ws.Range("B7").Value = 32
$ws.Range("C7").Value = "06:33"
$ws.Range("D7").Value = "11:07"
$ws.Range("E7").Value = "06:08:08"
$ws.Range("F7").Value = "06:10:25"
$ws.Range("G7").Value = "06:13:41"
$ws.Range("H7").Value = "06:16:58"
$ws.Range("I7").Value = "06:19:15"
$ws.Range("J7").Value = "2°"
$ws.Range("K7").Value = "06:08:33"
$ws.Range("L7").Value = -8
$ws.Range("M7").Value = "B"
$ws.Range("N7").Value = "3"

# ---------------------------------------------------------------------
# 3. Extend the conditional formatting to cover the new rows
# ---------------------------------------------------------------------
$ws.Range("L2:L5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("L2:L7"))
$ws.Range("N2:N5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("N2:N7"))
